$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "sobat_id"
$ws.Range("B1").Value = "posisi"
$ws.Range("C1").Value = "tahun_mitra"

# Data row
$ws.Range("A2").Value = "S003"
$ws.Range("B2").Value = "contoh999"
$ws.Range("C2").Value = 2025

# Match final selected cell as left by the author before saving
[void]$ws.Range("F5").Select()
